$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''65.528.34'
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("D3").Value = '''2.643.07'
$ws.Range("E3").Value = '  -1.18%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '''595.63'
$ws.Range("E5").Value = '  -0.89%  '
$ws.Range("D6").Value = '''155.56'
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '''0.630'
$ws.Range("E8").Value = '  +3.11%  '
$ws.Range("E9").Value = '  +2.97%  '
$ws.Range("D10").Value = '''5.81'
$ws.Range("E10").Value = '  -2.11%  '
$ws.Range("D11").Value = '''0.395'
$ws.Range("E11").Value = '  -1.26%  '
$ws.Range("E12").Value = '  +1.11%  '
$ws.Range("D13").Value = '''28.62'
$ws.Range("E13").Value = '  -3.45%  '
$ws.Range("D14").Value = '''0.0000197'
$ws.Range("E14").Value = '  -0.34%  '
$ws.Range("D15").Value = '''3.115.33'
$ws.Range("E15").Value = '  -1.15%  '
$ws.Range("D16").Value = '''65.361.52'
$ws.Range("E16").Value = '  -0.27%  '
$ws.Range("D17").Value = '''2.632.92'
$ws.Range("E17").Value = '  -1.52%  '
$ws.Range("D18").Value = '''12.54'
$ws.Range("E18").Value = '  -0.38%  '
$ws.Range("D19").Value = '''4.72'
$ws.Range("E19").Value = '  -2.42%  '
$ws.Range("D20").Value = '''7.41'
$ws.Range("E20").Value = '  -2.44%  '
$ws.Range("D21").Value = '''347.45'
$ws.Range("E21").Value = '  -1.41%  '
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").Value = '''68.80'
$ws.Range("E23").Value = '  -1.39%  '
$ws.Range("D24").Value = '''0.0000112'
$ws.Range("E24").Value = '  +0.15%  '
$ws.Range("E25").Value = '  -2.32%  '
$ws.Range("E26").Value = '  +2.51%  '
$ws.Range("E27").Value = '  -2.25%  '
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D28").Value = '''1.02'
$ws.Range("E28").Value = '  +1.85%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").Value = '''0.164'
$ws.Range("E29").Value = '  -2.55%  '
$ws.Range("D30").Value = '''7.86'
$ws.Range("E30").Value = '  -3.77%  '
$ws.Range("D31").Value = '''2.13'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("D32").Value = '''526.83'
$ws.Range("E32").Value = '  -3.18%  '
$ws.Range("E33").Value = '  -2.00%  '
$ws.Range("D34").Value = '''6.37'
$ws.Range("E34").Value = '  -3.80%  '
$ws.Range("D35").Value = '''5.39'
$ws.Range("E35").Value = '  -1.68%  '
$ws.Range("D37").Value = '''20.30'
$ws.Range("E37").Value = '  -0.87%  '
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '''1.92'
$ws.Range("E39").Value = '  -2.06%  '
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").Value = '''154.33'
$ws.Range("E40").Value = '  -2.27%  '
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").Value = '''160.21'
$ws.Range("E42").Value = '  -3.36%  '
$ws.Range("E43").Value = '  -1.03%  '
$ws.Range("D44").Value = '''0.0602'
$ws.Range("E44").Value = '  -2.29%  '
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("E46").Value = '  -3.67%  '
$ws.Range("E47").Value = '  -2.23%  '
$ws.Range("E48").Value = '  -2.90%  '
$ws.Range("D49").Value = '''0.0994'
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("E50").Value = '  +6.10%  '
$ws.Range("D51").Value = '''19.67'
$ws.Range("E51").Value = '  -1.91%  '
